$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C.
# This shifts the existing column C (the "rating detail" column) to column E,
# and leaves the existing column B (the "Jun_13" header / "UN" data column)
# untouched in place for now.
$ws.Range("C:D").Insert()

# The old B1 header ("Jun_13") needs to move into the new D1 position, since
# the two newly inserted weeks (Jun_17, Jun_15) become the new B1/C1 headers.
$ws.Range("D1").Value = $ws.Range("B1").Value()
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the two newly inserted columns (C and D) for every data row with the
# same "UN" placeholder value used throughout column B.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Give the two new columns (C, D) and the shifted column (E) the same ~8
# character width that column C originally had.
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
